# Update the cryptos list with newly scraped price / volume(1h) values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2;  D = "68.502.36"; E = "  +1.52%  " },
    @{ Row = 3;  D = "3.266.96";  E = $null },
    @{ Row = 4;  D = $null;       E = "  -0.01%  " },
    @{ Row = 5;  D = "583.88";    E = "  +0.91%  " },
    @{ Row = 6;  D = "181.72";    E = "  -0.34%  " },
    @{ Row = 8;  D = $null;       E = "  -0.06%  " },
    @{ Row = 9;  D = $null;       E = "  +0.67%  " },
    @{ Row = 10; D = "6.66";      E = "  -1.33%  " },
    @{ Row = 11; D = "0.424";     E = "  +2.09%  " },
    @{ Row = 12; D = $null;       E = "  +0.49%  " },
    @{ Row = 13; D = $null;       E = "  +0.20%  " },
    @{ Row = 14; D = "28.54";     E = "  +0.02%  " },
    @{ Row = 15; D = "68.515.13"; E = "  +1.51%  " },
    @{ Row = 16; D = "0.0000171"; E = "  +2.04%  " },
    @{ Row = 17; D = "3.196.11";  E = "  -1.55%  " },
    @{ Row = 18; D = "5.83";      E = "  -0.32%  " },
    @{ Row = 19; D = "13.55";     E = "  +0.01%  " },
    @{ Row = 20; D = "394.69";    E = "  +4.67%  " },
    @{ Row = 21; D = "7.69";      E = "  +0.83%  " },
    @{ Row = 22; D = "72.02";     E = "  +0.96%  " },
    @{ Row = 23; D = $null;       E = "  -0.06%  " },
    @{ Row = 24; D = "0.516";     E = "  +0.79%  " },
    @{ Row = 25; D = $null;       E = "  +0.47%  " },
    @{ Row = 26; D = "0.188";     E = "  +4.00%  " },
    @{ Row = 27; D = "9.63";      E = "  +0.16%  " },
    @{ Row = 28; D = $null;       E = "  -0.39%  " },
    @{ Row = 29; D = $null;       E = "  +0.03%  " },
    @{ Row = 30; D = "5.69";      E = "  -2.05%  " },
    @{ Row = 31; D = "22.94";     E = "  +1.01%  " },
    @{ Row = 32; D = "7.14";      E = "  +3.16%  " },
    @{ Row = 33; D = "1.28";      E = "  +0.34%  " },
    @{ Row = 34; D = $null;       E = "  +0.05%  " },
    @{ Row = 35; D = "164.40";    E = "  +0.53%  " },
    @{ Row = 36; D = $null;       E = "  +0.52%  " },
    @{ Row = 37; D = "1.91";      E = "  +2.58%  " },
    @{ Row = 38; D = "0.829";     E = "  -2.37%  " },
    @{ Row = 39; D = "4.60";      E = "  -0.46%  " },
    @{ Row = 40; D = "26.34";     E = "  -1.94%  " },
    @{ Row = 41; D = "6.54";      E = "  -4.52%  " },
    @{ Row = 42; D = "2.49";      E = "  -3.90%  " },
    @{ Row = 43; D = "41.30";     E = "  +1.03%  " },
    @{ Row = 44; D = $null;       E = "  +1.33%  " },
    @{ Row = 45; D = "346.35";    E = "  -3.28%  " },
    @{ Row = 46; D = "2.607.34";  E = "  -4.33%  " },
    @{ Row = 47; D = "24.62";     E = "  -3.43%  " },
    @{ Row = 48; D = "0.0281";    E = "  +0.41%  " },
    @{ Row = 49; D = "6.32";      E = "  +2.74%  " },
    @{ Row = 50; D = "31.58";     E = "  +0.98%  " },
    @{ Row = 51; D = $null;       E = "  -0.11%  " }
)

# Rows whose new Price text parses as a plain number. Those need the cell
# formatted as Text first, otherwise Excel (rightly) stores them as a
# number and the trailing-zero / "looks-like-a-float" text representation
# scraped from the site would be lost (e.g. "164.40" -> 164.4).
$numericLooking = @(5,6,10,11,14,16,18,19,20,21,22,24,26,27,30,31,32,33,35,37,38,39,40,41,42,43,45,47,48,49,50)

foreach ($u in $updates) {
    $r = $u.Row
    if ($null -ne $u.D) {
        $cell = $ws.Cells.Item($r, 4)
        if ($numericLooking -contains $r) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $u.D
    }
    if ($null -ne $u.E) {
        $ws.Cells.Item($r, 5).Value = $u.E
    }
}
